$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# ------------------------------------------------------------------
# Insert two brand-new slides right after the title slide (position 2
# and 3). This pushes the existing "Objectifs du projet" slide down to
# position 4 while it keeps its original SlideID (257).
# ------------------------------------------------------------------
$sEquipe = $p.Slides.AddSlide(2, $layout)
$sMission = $p.Slides.AddSlide(3, $layout)

# ====================================================================
# Slide 2 : "L'Equipe"
# ====================================================================
$titre = $sEquipe.Shapes.Item(1)
$titre.Name = "Titre 1"
$titre.TextFrame.TextRange.Text = "L’Équipe"
$titrePara = $titre.TextFrame.TextRange.Paragraphs(1,1)
$titrePara.ParagraphFormat.Alignment = 2

$corps = $sEquipe.Shapes.Item(2)
$corps.Name = "Espace réservé du contenu 2"
$corps.TextFrame.TextRange.Text = "Chef de projet: `rAlban INQUEL`r`rMembre de l’équipe :`rChristie BUNLON`rNicolas CHATELAIN`rBrice HOFFMANN`rAlban INQUEL"
$ctr = $corps.TextFrame.TextRange
$ctr.Paragraphs(2,1).IndentLevel = 1
$ctr.Paragraphs(6,1).IndentLevel = 1
$ctr.Paragraphs(7,1).IndentLevel = 1
$ctr.Paragraphs(8,1).IndentLevel = 1
$ctr.Paragraphs(9,1).IndentLevel = 1

# ====================================================================
# Slide 3 : "Misssion"
# ====================================================================
$titre2 = $sMission.Shapes.Item(1)
$titre2.Name = "Titre 1"
$titre2.TextFrame.TextRange.Text = "Misssion"
$titre2Para = $titre2.TextFrame.TextRange.Paragraphs(1,1)
$titre2Para.ParagraphFormat.Alignment = 2

$corps2 = $sMission.Shapes.Item(2)
$corps2.Name = "Espace réservé du contenu 2"
$corps2.TextFrame.TextRange.Text = "`r`r`rInterception de datagrammes destinés à une machine `ret affichage de leur contenu en temps réel.`r"
$ctr2 = $corps2.TextFrame.TextRange
$ctr2.Paragraphs(4,1).ParagraphFormat.Alignment = 2
$p5 = $ctr2.Paragraphs(5,1)
$p5.ParagraphFormat.Alignment = 2
$p5.ParagraphFormat.Bullet.Visible = 0

# ====================================================================
# Slide 4 (formerly slide 2) : "Objectifs" (was "Objectifs du projet")
# ====================================================================
$sObjectifs = $p.Slides.Item(4)

$titre3 = $sObjectifs.Shapes.Item(1)
$titre3.TextFrame.TextRange.Text = "Objectifs"
$titre3Para = $titre3.TextFrame.TextRange.Paragraphs(1,1)
$titre3Para.ParagraphFormat.Alignment = 2

$corps3 = $sObjectifs.Shapes.Item(2)
$corps3.TextFrame.AutoSize = 2
$lines = @(
    "Pouvoir scanner l’intégralité du réseau avant le 15/04.",
    "Pouvoir modifier la table ARP à distance avant le 22/04.",
    "Pouvoir intercepter et rerouter les paquets avant le 29/04.",
    "Création d’une interface graphique avant le 03/05.",
    "Pouvoir obtenir les mots de passe transitant en clair avant le 06/05.",
    "Pouvoir modifier les réponses DNS avant le 15/05.",
    "Pouvoir intégrer les cookies de la « victime » automatiquement sur la machine attaquante avant le 23/05.",
    "Affichage des pages web visitées avant le 30/05.",
    "Pouvoir intercepter les transmissions VoIP avant le 05/06.",
    "Pouvoir générer de faux certificats automatiquement avant le 13/06."
)
$corps3.TextFrame.TextRange.Text = [string]::Join("`r", $lines)
$ctr3 = $corps3.TextFrame.TextRange
for ($i = 1; $i -le $lines.Count; $i++) {
    $ctr3.Paragraphs($i,1).IndentLevel = 0
}

# ====================================================================
# Slides 5-9 : blank placeholder slides appended at the end of the deck
# ====================================================================
$p.Slides.AddSlide(5, $layout) | Out-Null
$p.Slides.AddSlide(6, $layout) | Out-Null
$p.Slides.AddSlide(7, $layout) | Out-Null
$p.Slides.AddSlide(8, $layout) | Out-Null
$p.Slides.AddSlide(9, $layout) | Out-Null

$sBlank5 = $p.Slides.Item(5)
$sBlank5.Shapes.Item(1).Name = "Titre 1"
$sBlank5.Shapes.Item(2).Name = "Espace réservé du contenu 2"
$sBlank5.Shapes.Item(1).Left = 656502 / 12700.0
$sBlank5.Shapes.Item(1).Top = 442327 / 12700.0
$sBlank5.Shapes.Item(1).Width = 9404723 / 12700.0
$sBlank5.Shapes.Item(1).Height = 1400530 / 12700.0

for ($i = 6; $i -le 9; $i++) {
    $sb = $p.Slides.Item($i)
    $sb.Shapes.Item(1).Name = "Titre 1"
    $sb.Shapes.Item(2).Name = "Espace réservé du contenu 2"
}

Write-Output "Done. Slide count: $($p.Slides.Count)"
